# Regenerate the "example output" worksheet in place.
#
# The source script that produced this workbook was re-run against a
# newer extract: it now emits an extra "Matching Filter" column (right
# after "Distance from track (km)") and two extra POI rows ("Les
# Acacias" and "L'Islette - Fondettes"), while the header row lost its
# ad-hoc bold/bordered/centered style (the sheet is emitted as plain
# data now, no manual formatting).
#
# Rather than trying to replay insert/shift operations cell-by-cell,
# just clear the sheet (content AND formatting) and write the final
# grid directly - this is effectively what re-running the export
# script produced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe all existing cell values/styles on the sheet.
$ws.Cells.Clear()

# ---- Header row (row 1) - plain, no style ----
$ws.Cells.Item(1, 1).Value = 'Kilometers from start'
$ws.Cells.Item(1, 2).Value = 'Distance from track (km)'
$ws.Cells.Item(1, 3).Value = 'Matching Filter'
$ws.Cells.Item(1, 4).Value = 'Name'
$ws.Cells.Item(1, 5).Value = 'Website'
$ws.Cells.Item(1, 6).Value = 'Phone'
$ws.Cells.Item(1, 7).Value = 'Opening hours'
$ws.Cells.Item(1, 8).Value = 'OSM Tags'
$ws.Cells.Item(1, 9).Value = 'lat'
$ws.Cells.Item(1, 10).Value = 'lon'

# ---- Row 2: Les Acacias (new row) ----
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = 4.62
$ws.Cells.Item(2, 3).Value = 'tourism=camp_site'
$ws.Cells.Item(2, 4).Value = 'Les Acacias'
$ws.Cells.Item(2, 5).Value = 'https://www.camping-tours.fr/'
$ws.Cells.Item(2, 6).Value = '+33 2 47 44 08 16'
$ws.Cells.Item(2, 7).Value = ""
$ws.Cells.Item(2, 8).Value = '{''addr:city'': ''La Ville-aux-Dames'', ''addr:postcode'': ''37700'', ''addr:street'': ''Rue Berthe Morisot'', ''barrier'': ''fence'', ''caravans'': ''yes'', ''email'': ''contact@camplvad.com'', ''internet_access'': ''yes'', ''internet_access:fee'': ''no'', ''name'': ''Les Acacias'', ''phone'': ''+33 2 47 44 08 16'', ''stars'': ''3'', ''tents'': ''yes'', ''tourism'': ''camp_site'', ''website'': ''https://www.camping-tours.fr/''}'
$ws.Cells.Item(2, 9).Value = 47.4020858
$ws.Cells.Item(2, 10).Value = 0.7801299

# ---- Row 3: Aire camping-car de Saint Avertin (was row 2) ----
$ws.Cells.Item(3, 1).Value = 1.57
$ws.Cells.Item(3, 2).Value = 1.57
$ws.Cells.Item(3, 3).Value = 'tourism=camp_site'
$ws.Cells.Item(3, 4).Value = 'Aire camping-car de Saint Avertin'
$ws.Cells.Item(3, 5).Value = 'https://www.onlypark.fr/aire-camping-car-de-st-avertin/'
$ws.Cells.Item(3, 6).Value = '+33 2 47 27 87 47'
$ws.Cells.Item(3, 7).Value = ""
$ws.Cells.Item(3, 8).Value = '{''addr:city'': ''Saint-Avertin'', ''capacity'': ''20'', ''caravans'': ''yes'', ''charge'': ''19 EUR'', ''charge:conditional'': ''12 EUR @ (Sep-Jun)'', ''contact:email'': ''campingtoursvaldeloire@onlycamp.fr'', ''contact:phone'': ''+33 2 47 27 87 47'', ''contact:website'': ''https://www.onlypark.fr/aire-camping-car-de-st-avertin/'', ''drinking_water'': ''yes'', ''fee'': ''yes'', ''internet_access'': ''wlan'', ''name'': ''Aire camping-car de Saint Avertin'', ''network'': ''Onlypark'', ''operator'': ''onlycamp'', ''power_supply'': ''yes'', ''power_supply:charge'': ''3 EUR/4 hours'', ''power_supply:fee'': ''yes'', ''power_supply:maxcurrent'': ''10'', ''sanitary_dump_station'': ''yes'', ''sanitary_dump_station:charge'': ''3 EUR/20 minutes'', ''sanitary_dump_station:fee'': ''yes'', ''shower'': ''yes'', ''stars'': ''4'', ''tents'': ''yes'', ''toilets'': ''no'', ''tourism'': ''camp_site'', ''water_point'': ''yes'', ''wheelchair'': ''yes''}'
$ws.Cells.Item(3, 9).Value = 47.3708862
$ws.Cells.Item(3, 10).Value = 0.7243202

# ---- Row 4: L'Islette - Fondettes (new row) ----
$ws.Cells.Item(4, 1).Value = 9.91
$ws.Cells.Item(4, 2).Value = 4.44
$ws.Cells.Item(4, 3).Value = 'tourism=camp_site'
$ws.Cells.Item(4, 4).Value = 'L''Islette - Fondettes'
$ws.Cells.Item(4, 5).Value = ""
$ws.Cells.Item(4, 6).Value = ""
$ws.Cells.Item(4, 7).Value = ""
$ws.Cells.Item(4, 8).Value = '{''name'': "L''Islette - Fondettes", ''tourism'': ''camp_site''}'
$ws.Cells.Item(4, 9).Value = 47.3892756
$ws.Cells.Item(4, 10).Value = 0.5959279

# ---- Row 5: Camping La Mignardière (was row 3) ----
$ws.Cells.Item(5, 1).Value = 11.11
$ws.Cells.Item(5, 2).Value = 0.08
$ws.Cells.Item(5, 3).Value = 'tourism=camp_site'
$ws.Cells.Item(5, 4).Value = 'Camping La Mignardière'
$ws.Cells.Item(5, 5).Value = 'https://www.mignardiere.com/'
$ws.Cells.Item(5, 6).Value = '+33 2 47 73 31 00'
$ws.Cells.Item(5, 7).Value = ""
$ws.Cells.Item(5, 8).Value = '{''addr:city'': ''Ballan-Miré'', ''addr:housenumber'': ''22'', ''addr:postcode'': ''37510'', ''addr:street'': ''Avenue des Aubépines'', ''cabins'': ''yes'', ''capacity:caravans'': ''114'', ''capacity:tents'': ''114'', ''caravans'': ''yes'', ''drinking_water'': ''yes'', ''motorhome'': ''yes'', ''name'': ''Camping La Mignardière'', ''phone'': ''+33 2 47 73 31 00'', ''sanitary_dump_station'': ''yes'', ''shower'': ''yes'', ''stars'': ''4'', ''tents'': ''yes'', ''toilets'': ''yes'', ''tourism'': ''camp_site'', ''washing_machine'': ''yes'', ''website'': ''https://www.mignardiere.com/''}'
$ws.Cells.Item(5, 9).Value = 47.3557614
$ws.Cells.Item(5, 10).Value = 0.6332265
